# Updated cryptos list on Thu Feb 29 03:26:03 UTC 2024 with GitHub Actions
#
# Refreshes the Price (D) / Volume(1h) (E) columns of the crypto table, and
# re-orders three rows (Toncoin moves above Cosmos/Dai) to reflect the new
# ranking, exactly as produced by the upstream data-refresh job.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Columns D/E sometimes hold strings that *look* numeric (e.g. "412.38",
    # "0.576", "  +4.14%  "). Writing straight to .Value lets the COM layer
    # infer a real number/percentage and silently reformat it, so force the
    # cell to Text first, assign, then drop the format back to the sheet's
    # normal (General) style so no stray formatting is left behind.
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = "Normal"
}

# row -> @{ column letter = new value }
$changes = [ordered]@{
    2  = @{ D = "61.395.28"; E = "  +7.83%  " }
    3  = @{ D = "3.405.65";  E = "  +4.96%  " }
    4  = @{ E = "  +0.13%  " }
    5  = @{ D = "412.38";    E = "  +4.14%  " }
    6  = @{ D = "122.22";    E = "  +13.43%  " }
    7  = @{ D = "3.399.71";  E = "  +4.91%  " }
    8  = @{ D = "0.576";     E = "  -0.93%  " }
    9  = @{ E = "  +0.11%  " }
    10 = @{ D = "0.638";     E = "  +3.24%  " }
    11 = @{ D = "0.115";     E = "  +20.77%  " }
    12 = @{ D = "41.05";     E = "  +4.82%  " }
    13 = @{ E = "  -0.66%  " }
    14 = @{ D = "3.952.86";  E = "  +5.46%  " }
    15 = @{ D = "8.38";      E = "  +1.15%  " }
    16 = @{ D = "19.46";     E = "  +3.20%  " }
    17 = @{ D = "3.409.94";  E = "  +5.08%  " }
    18 = @{ D = "61.406.58"; E = "  +8.23%  " }
    19 = @{ E = "  -0.61%  " }
    20 = @{ D = "10.81";     E = "  -1.09%  " }
    21 = @{ D = "0.0000120"; E = "  +9.27%  " }
    22 = @{ E = "  -0.18%  " }
    23 = @{ D = "12.78";     E = "  -0.23%  " }
    24 = @{ D = "297.95";    E = "  +1.99%  " }
    25 = @{ D = "76.12";     E = "  +2.63%  " }
    26 = @{ D = "3.14";      E = "  -1.08%  " }
    27 = @{ D = "30.64";     E = "  +9.34%  " }
    28 = @{ E = "  +12.07%  " }
    29 = @{ D = "4.27";      E = "  -2.39%  " }
    30 = @{ D = "7.63";      E = "  -4.30%  " }
    31 = @{ D = "0.172";     E = "  +1.72%  " }
    32 = @{ D = "0.116";     E = "  +5.58%  " }
    33 = @{ D = "42.43";     E = "  +3.44%  " }
    34 = @{ B = "Toncoin"; C = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"; D = "2.55";  E = "  +19.69%  " }
    35 = @{ B = "Cosmos";  C = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"; D = "11.38"; E = "  +2.04%  " }
    36 = @{ B = "Dai";     C = "https://coinranking.com/coin/MoTuySvg7+dai-dai"; D = "1.00"; E = "  +0.03%  " }
    37 = @{ D = "0.0479";   E = "  -0.64%  " }
    38 = @{ D = "52.48";    E = "  +2.70%  " }
    39 = @{ D = "3.53";     E = "  +2.23%  " }
    40 = @{ E = "  -0.01%  " }
    41 = @{ D = "3.00";     E = "  +1.11%  " }
    42 = @{ E = "  +5.32%  " }
    43 = @{ D = "0.122";    E = "  +0.76%  " }
    44 = @{ D = "133.37";   E = "  -2.79%  " }
    45 = @{ D = "17.16";    E = "  +3.15%  " }
    46 = @{ D = "3.92";     E = "  +0.08%  " }
    47 = @{ D = "0.282";    E = "  -0.19%  " }
    48 = @{ D = "2.19";     E = "  -1.69%  " }
    49 = @{ D = "21.81";    E = "  -2.13%  " }
    50 = @{ D = "2.203.74"; E = "  +2.34%  " }
    51 = @{ D = "3.747.85"; E = "  +5.22%  " }
}

foreach ($row in $changes.Keys) {
    $cols = $changes[$row]
    foreach ($col in $cols.Keys) {
        $addr = "$col$row"
        if ($col -eq "D") {
            Set-TextValue $ws.Range($addr) $cols[$col]
        } else {
            $ws.Range($addr).Value = $cols[$col]
        }
    }
}
